$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above row 355, pushing the existing rows
# (355-362) down to (357-364) and preserving their data untouched.
$ws.Rows.Item(355).Insert()
$ws.Rows.Item(355).Insert()

# New row 355: Apio, Primera, week of 2021-09-09 (serial 44448)
$ws.Cells.Item(355,1).Value  = 6
$ws.Cells.Item(355,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(355,3).Value  = "Metropolitana"
$ws.Cells.Item(355,4).Value  = 44448
$ws.Cells.Item(355,5).Value  = 13
$ws.Cells.Item(355,6).Value  = 100112017
$ws.Cells.Item(355,7).Value  = "Apio"
$ws.Cells.Item(355,8).Value  = "Americana (o)"
$ws.Cells.Item(355,9).Value  = "Primera"
$ws.Cells.Item(355,10).Value = 3000
$ws.Cells.Item(355,11).Value = 7000
$ws.Cells.Item(355,12).Value = 8000
$ws.Cells.Item(355,13).Value = 7600
$ws.Cells.Item(355,14).Value = "$/docena de matas"
$ws.Cells.Item(355,15).Value = "Región de Coquimbo"
$ws.Cells.Item(355,16).Value = 1267
$ws.Cells.Item(355,17).Value = 6
$ws.Cells.Item(355,18).Value = "Hortaliza"

# New row 356: Apio, Segunda, week of 2021-09-09 (serial 44448)
$ws.Cells.Item(356,1).Value  = 6
$ws.Cells.Item(356,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(356,3).Value  = "Metropolitana"
$ws.Cells.Item(356,4).Value  = 44448
$ws.Cells.Item(356,5).Value  = 13
$ws.Cells.Item(356,6).Value  = 100112017
$ws.Cells.Item(356,7).Value  = "Apio"
$ws.Cells.Item(356,8).Value  = "Americana (o)"
$ws.Cells.Item(356,9).Value  = "Segunda"
$ws.Cells.Item(356,10).Value = 1200
$ws.Cells.Item(356,11).Value = 6000
$ws.Cells.Item(356,12).Value = 6000
$ws.Cells.Item(356,13).Value = 6000
$ws.Cells.Item(356,14).Value = "$/docena de matas"
$ws.Cells.Item(356,15).Value = "Región de Coquimbo"
$ws.Cells.Item(356,16).Value = 1000
$ws.Cells.Item(356,17).Value = 6
$ws.Cells.Item(356,18).Value = "Hortaliza"
